$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing data row down to the new rows first
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E8").PasteSpecial(-4122)

# New row 7
$ws.Range("A7").Value = 21357
$ws.Range("B7").Value = "HIJRI"
$ws.Range("C7").Value = "LENA"
$ws.Range("D7").Value = 123
$ws.Range("E7").Value = "XI-MM-2"

# New row 8
$ws.Range("A8").Value = 213
$ws.Range("B8").Value = 232
$ws.Range("C8").Value = "NIH"
$ws.Range("D8").Value = 123
$ws.Range("E8").Value = "X-LPB-2"

# Column width adjustments: D and E become a single uniform width (target stored width 16.63)
$ws.Range("D:E").ColumnWidth = 15.8
